$d = $word.ActiveDocument

# 1) Replace the old sentence with the new full sentence. Using $d.Content as
#    the Find range means that after Execute() it collapses to exactly the
#    span of the replacement text (Start/End), so we can locate it precisely
#    without fragile text comparisons (Paragraph.Range.Text carries a
#    trailing paragraph-mark char that complicates equality checks).
$rng = $d.Content
$old = "En una localidad viven muchas personas"
$new = "En una localidad pueden (0…1) vivir muchas personas"
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

$pStart = $rng.Start

# 2) The final text is split (per the source edit) into five runs that all
#    share identical run formatting:
#      "En una localidad " | "pueden (0…1) " | "viv" | "ir" | " muchas personas"
#    The engine coalesces freshly-typed adjacent same-format text back into a
#    single run, so force the boundaries to "stick" by toggling a formatting
#    property on and back off across each split point (net effect on the
#    formatting is a no-op, but it breaks the run apart for good).
$seg1 = "En una localidad "
$seg2 = "pueden (0…1) "
$seg3 = "viv"
$seg4 = "ir"
$seg5 = " muchas personas"

$p1 = $pStart + $seg1.Length
$p2 = $p1 + $seg2.Length
$p3 = $p2 + $seg3.Length
$p4 = $p3 + $seg4.Length

$b1 = $d.Range($pStart, $p1)
$b1.Font.Bold = 1
$b1.Font.Bold = 0

$b2 = $d.Range($p1, $p2)
$b2.Font.Bold = 1
$b2.Font.Bold = 0

$b3 = $d.Range($p2, $p3)
$b3.Font.Bold = 1
$b3.Font.Bold = 0

$b4 = $d.Range($p3, $p4)
$b4.Font.Bold = 1
$b4.Font.Bold = 0
